$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.874784666666667
$ws.Range("N2").Value = 8.624354
$ws.Range("O2").Value = 0.1187109652550681
$ws.Range("P2").Value = 0.121184727686443
$ws.Range("Q2").Value = 1.116570197579556
$ws.Range("R2").Value = 10.049131778216
$ws.Range("S2").Value = 0.1187109652550681
$ws.Range("T2").Value = 0.121184727686443

# Row 3
$ws.Range("O3").Value = 0.4442422727481699
$ws.Range("P3").Value = 0.4534996302499962
$ws.Range("S3").Value = 0.4442422727481699
$ws.Range("T3").Value = 0.4534996302499962

# Row 4
$ws.Range("M4").Value = 5.147441999999999
$ws.Range("N4").Value = 15.442326
$ws.Range("O4").Value = 0.2125577666737049
$ws.Range("P4").Value = 0.2169871588243338
$ws.Range("Q4").Value = 1.999273336056
$ws.Range("R4").Value = 17.993460024504
$ws.Range("S4").Value = 0.2125577666737049
$ws.Range("T4").Value = 0.2169871588243338

# Row 5
$ws.Range("M5").Value = 1.483016
$ws.Range("N5").Value = 2.966032
$ws.Range("O5").Value = 0.06123946008548931
$ws.Range("P5").Value = 0.04167706708575228
$ws.Range("Q5").Value = 0.5760053917546668
$ws.Range("R5").Value = 3.456032350528001
$ws.Range("S5").Value = 0.06123946008548931
$ws.Range("T5").Value = 0.04167706708575228

# Row 6
$ws.Range("M6").Value = 3.953360666666667
$ws.Range("N6").Value = 11.860082
$ws.Range("O6").Value = 0.1632495352375677
$ws.Range("P6").Value = 0.1666514161534747
$ws.Range("Q6").Value = 1.535490554080889
$ws.Range("R6").Value = 13.819414986728
$ws.Range("S6").Value = 0.1632495352375677
$ws.Range("T6").Value = 0.1666514161534747
